# Helper: locate a shape's 1-based collection index by its stable .Id
function Find-ShapeIndexById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        if ($shapes.Item($i).Id -eq $id) {
            return $i
        }
    }
    return -1
}

# Helper: locate a shape object by its stable .Id
function Find-ShapeById($shapes, $id) {
    $idx = Find-ShapeIndexById $shapes $id
    if ($idx -eq -1) {
        return $null
    }
    return $shapes.Item($idx)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1. Shrink the "AWS Cloud" bounding rectangle (id 25) to fit the tightened layout.
#    (target size 8034612 x 3925693 EMU)
$rect = Find-ShapeById $s.Shapes 25
$rect.Width = 632.6466369732284
$rect.Height = 309.1096954393701

# 2. Remove the "Amazon EC2" icon (id 36) and its caption textbox (id 37).
$ec2Icon = Find-ShapeById $s.Shapes 36
$ec2Icon.Delete()
$ec2Label = Find-ShapeById $s.Shapes 37
$ec2Label.Delete()

# 3. Reposition the "Graphic 17" icon (id 44) and its "TextBox 11" caption (id 45)
#    down into the second row of the diagram.
#    (target off 7234289,2914580 / 6469114,3658548 EMU)
$graphic17 = Find-ShapeById $s.Shapes 44
$graphic17.Left = 569.6290589181102
$graphic17.Top = 229.4944881889764

$textBox45 = Find-ShapeById $s.Shapes 45
$textBox45.Left = 509.3790551181102
$textBox45.Top = 288.07464566929133

# 4. Group the "AWS IAM" icon (id 30) with its caption (id 31) into "Group 9"
#    and move the resulting group into place. (target off 1297526,2910805 EMU)
$idxA = Find-ShapeIndexById $s.Shapes 30
$idxB = Find-ShapeIndexById $s.Shapes 31
$group9 = $s.Shapes.Range(@($idxA, $idxB)).Group()
$group9.Name = "Group 9"
$group9.Left = 102.16740157480315
$group9.Top = 229.19724409448818

# 5. Group the "AWS Systems Manager" icon (id 34) with its caption (id 35) into
#    "Group 11" and move the resulting group into place.
#    (target off 4759829,2909075 EMU)
$idxA = Find-ShapeIndexById $s.Shapes 34
$idxB = Find-ShapeIndexById $s.Shapes 35
$group11 = $s.Shapes.Range(@($idxA, $idxB)).Group()
$group11.Name = "Group 11"
$group11.Left = 374.78968813937007
$group11.Top = 229.06102762204722

# 6. Group the "Amazon EventBridge" icon (id 40) with its caption (id 41) into
#    "Group 10" and move the resulting group into place.
#    (target off 2834433,2910805 EMU)
$idxA = Find-ShapeIndexById $s.Shapes 40
$idxB = Find-ShapeIndexById $s.Shapes 41
$group10 = $s.Shapes.Range(@($idxA, $idxB)).Group()
$group10.Name = "Group 10"
$group10.Left = 223.18370828740157
$group10.Top = 229.19724409448818
